# Cambio de clave de priorizacion en offset
#
# Updates the "Maquinas" sheet: a new "Descartonadora 3" row is inserted
# (logically) after "Descartonadora 2", pushing the old "Ventana" /
# "Pegado" rows down by one, and a brand-new "Pegado" row (with new
# values) is appended at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Maquinas")
$ws.Activate()

# --- Row 3: Impresion Offset ---
$ws.Cells.Item(3, 3).Value = 3500
$ws.Cells.Item(3, 4).Value = 60
$ws.Cells.Item(3, 5).Value = 45

# --- Row 4: Impresion Flexo ---
$ws.Cells.Item(4, 3).Value = 1100
$ws.Cells.Item(4, 4).Value = 45

# --- Row 6: Troquelado / Automatica ---
$ws.Cells.Item(6, 3).Value = 4000

# --- Row 7: Troquelado / Manual 1 ---
$ws.Cells.Item(7, 3).Value = 750
$ws.Cells.Item(7, 4).Value = 90
$ws.Cells.Item(7, 5).Value = 40

# --- Row 8: Troquelado / Manual 2 ---
$ws.Cells.Item(8, 3).Value = 750
$ws.Cells.Item(8, 4).Value = 90
$ws.Cells.Item(8, 5).Value = 40

# --- Row 9: Troquelado / Manual 3 ---
$ws.Cells.Item(9, 3).Value = 750
$ws.Cells.Item(9, 4).Value = 90
$ws.Cells.Item(9, 5).Value = 40

# --- Row 12: was Ventana/Ventanas, now Descartonado/Descartonadora 3 ---
$ws.Cells.Item(12, 1).Value = "Descartonado"
$ws.Cells.Item(12, 2).Value = "Descartonadora 3"
$ws.Cells.Item(12, 3).Value = 3000
$ws.Cells.Item(12, 4).Value = 15
$ws.Cells.Item(12, 5).Value = 7

# --- Row 13: was Pegado/Pegadora 1, now Ventana/Ventanas ---
$ws.Cells.Item(13, 1).Value = "Ventana"
$ws.Cells.Item(13, 2).Value = "Ventanas"
$ws.Cells.Item(13, 3).Value = 1500
$ws.Cells.Item(13, 4).Value = 25
$ws.Cells.Item(13, 5).Value = 10

# --- Row 14 (new): Pegado/Pegadora 1 ---
$ws.Cells.Item(14, 1).Value = "Pegado"
$ws.Cells.Item(14, 2).Value = "Pegadora 1"
$ws.Cells.Item(14, 3).Value = 4000
$ws.Cells.Item(14, 4).Value = 20
$ws.Cells.Item(14, 5).Value = 8

# Update the active selection left behind on this sheet (was D15, now C14).
$ws.Range("C14").Select()
